# Improve test coverage for xlsx read
#
# - Re-applies the "Normal" style to the header + data block (A1:E3) so
#   those cells pick up a distinct style record (exercises style reuse
#   on read).
# - Gives D3 a custom number format (negative numbers in red) to exercise
#   custom numFmt parsing.
# - Merges A3:B3 (clearing B3's value, same as Excel does on merge) to
#   exercise merged-cell reading.
# - Adds a new row 4 with a date cell and a date+time cell to exercise
#   date/datetime parsing.
# - Widens column B a bit and moves the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reapply the default/"Normal" style across the existing header + data rows.
$ws.Range("A1:E3").Style = "Normal"

# Custom numeric format (negative values shown in red) on D3.
$ws.Range("D3").NumberFormat = "#,##0.00;[RED]-#,##0.00"

# Merge A3:B3 (Excel keeps only the upper-left value, clearing B3).
$ws.Range("A3:B3").Merge()

# New row of date / datetime test values.
$ws.Range("A4").Value = 43955
$ws.Range("A4").NumberFormat = "yyyy-mm-dd"

$ws.Range("B4").Value = 43955.5626388889
$ws.Range("B4").NumberFormat = "yyyy-mm-dd hh:mm:ss"

# Widen column B slightly.
$ws.Columns.Item(2).ColumnWidth = 16.95

# Move the active selection to A3.
$ws.Range("A3").Select() | Out-Null
